# Update countries & provincias Spain
# Refresh case counts for several countries and update the "last updated"
# timestamp. Because the sheet is kept sorted by "Casos totales" (column B)
# descending, two pairs of neighbouring rows swap places once their counts
# are refreshed (Kazajistan/Uzbekistan and Tunez/Bulgaria).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A4:A216")

function Set-CountryRow {
    param($CountryName, $B, $C, $D, $E, $F, $G, $H)
    $cell = $dataRange.Find($CountryName)
    $r = $cell.Row()
    $ws.Cells.Item($r, 2).Value = $B
    $ws.Cells.Item($r, 3).Value = $C
    $ws.Cells.Item($r, 4).Value = $D
    $ws.Cells.Item($r, 5).Value = $E
    $ws.Cells.Item($r, 6).Value = $F
    $ws.Cells.Item($r, 7).Value = $G
    $ws.Cells.Item($r, 8).Value = $H
}

# Straightforward refreshes (row position unchanged)
Set-CountryRow "Estados Unidos" 644348 259 48708 567086 13487 25 28554
Set-CountryRow "Israel"         12501  0   2563  9806   180   2  132
Set-CountryRow "India"          12456  86  1513  10520  0     1  423
Set-CountryRow "Hungria"        1652   73  199   1311   58    8  142

# Kazajistan / Uzbekistan: Uzbekistan's refreshed count overtakes Kazajistan's,
# so it now sorts above it. Update the country label together with the stats
# on both rows so the table stays sorted by column B descending.
$kzCell = $dataRange.Find("Kazajistan")
$kzRow = $kzCell.Row()
$uzCell = $dataRange.Find("Uzbekistan")
$uzRow = $uzCell.Row()

$ws.Cells.Item($kzRow, 1).Value = "Uzbekistan"
$ws.Cells.Item($kzRow, 2).Value = 1349
$ws.Cells.Item($kzRow, 3).Value = 47
$ws.Cells.Item($kzRow, 4).Value = 107
$ws.Cells.Item($kzRow, 5).Value = 1238
$ws.Cells.Item($kzRow, 6).Value = 8
$ws.Cells.Item($kzRow, 7).Value = 0
$ws.Cells.Item($kzRow, 8).Value = 4

$ws.Cells.Item($uzRow, 1).Value = "Kazajistan"
$ws.Cells.Item($uzRow, 2).Value = 1331
$ws.Cells.Item($uzRow, 3).Value = 36
$ws.Cells.Item($uzRow, 4).Value = 240
$ws.Cells.Item($uzRow, 5).Value = 1075
$ws.Cells.Item($uzRow, 6).Value = 22
$ws.Cells.Item($uzRow, 7).Value = 0
$ws.Cells.Item($uzRow, 8).Value = 16

# Tunez / Bulgaria: Bulgaria's refreshed count overtakes Tunez's, so it now
# sorts above it.
$tnCell = $dataRange.Find("Tunez")
$tnRow = $tnCell.Row()
$bgCell = $dataRange.Find("Bulgaria")
$bgRow = $bgCell.Row()

$ws.Cells.Item($tnRow, 1).Value = "Bulgaria"
$ws.Cells.Item($tnRow, 2).Value = 783
$ws.Cells.Item($tnRow, 3).Value = 36
$ws.Cells.Item($tnRow, 4).Value = 122
$ws.Cells.Item($tnRow, 5).Value = 624
$ws.Cells.Item($tnRow, 6).Value = 31
$ws.Cells.Item($tnRow, 7).Value = 1
$ws.Cells.Item($tnRow, 8).Value = 37

$ws.Cells.Item($bgRow, 1).Value = "Tunez"
$ws.Cells.Item($bgRow, 2).Value = 780
$ws.Cells.Item($bgRow, 3).Value = 0
$ws.Cells.Item($bgRow, 4).Value = 43
$ws.Cells.Item($bgRow, 5).Value = 702
$ws.Cells.Item($bgRow, 6).Value = 89
$ws.Cells.Item($bgRow, 7).Value = 0
$ws.Cells.Item($bgRow, 8).Value = 35

# Update the "last refreshed" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Abril de 2020 a las 07:52"
